$wb = $excel.ActiveWorkbook

# Rename the first two sheets
$wb.Worksheets.Item(1).Name = "Cristi's tests (12May08)"
$wb.Worksheets.Item(2).Name = "Corneliu's changes (16May08)"

# Select the second sheet and change the active cell / selection
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

# Re-merge these header cells so they move to the end of the mergeCells list,
# in the order G, H, K, L, M, J, I (matching the target OOXML ordering).
$ws2.Range("G2:G3").UnMerge()
$ws2.Range("H2:H3").UnMerge()
$ws2.Range("K2:K3").UnMerge()
$ws2.Range("L2:L3").UnMerge()
$ws2.Range("M2:M3").UnMerge()
$ws2.Range("J2:J3").UnMerge()
$ws2.Range("I2:I3").UnMerge()

$ws2.Range("G2:G3").Merge()
$ws2.Range("H2:H3").Merge()
$ws2.Range("K2:K3").Merge()
$ws2.Range("L2:L3").Merge()
$ws2.Range("M2:M3").Merge()
$ws2.Range("J2:J3").Merge()
$ws2.Range("I2:I3").Merge()

$ws2.Range("H11").Select()
